$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the "Notes" column from E to G (shifted right to make room for new
# Species / Life stage columns in the MEASUREMENTS section). Copy() carries
# over both the shared-string value and the cell style/format.
$ws.Range("E8").Copy($ws.Range("G8"))
$ws.Range("E12").Copy($ws.Range("G12"))
$ws.Range("E16").Copy($ws.Range("G16"))
$ws.Range("E17").Copy($ws.Range("G17"))

$ws.Range("E8").Clear()
$ws.Range("E12").Clear()
$ws.Range("E16").Clear()

# Add new Species / Life stage columns to the MEASUREMENTS section header
# (row 17): copy the formatting from the neighboring header cell (C17, bold)
# then set the new text.
$ws.Range("C17").Copy($ws.Range("D17"))
$ws.Range("C17").Copy($ws.Range("E17"))
$ws.Range("D17").Value2 = "Species"
$ws.Range("E17").Value2 = "Life stage"

# ... and the corresponding data row (row 18), copying formatting from C18.
# Write E18 before D18 so the shared-string table order matches the source
# authoring order (Life stage, Adult/smolt, str).
$ws.Range("C18").Copy($ws.Range("D18"))
$ws.Range("C18").Copy($ws.Range("E18"))
$ws.Range("E18").Value2 = "Adult/smolt"
$ws.Range("D18").Value2 = "str"

$ws.Range("G8").Select()
